# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# worksheets, which hold identical data tables.

$wb = $excel.ActiveWorkbook

# Row number (key) -> new F-column value (value)
$updates = @{
    2 = 370
    3 = 1243
    4 = 1481
    5 = 56
    6 = 6132
    7 = 104
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
